$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9391857.142857142
$ws.Range("C3").Value = 13571.42857142864
$ws.Range("C4").Value = 3582500
$ws.Range("C6").Value = 6232142.857142857
$ws.Range("C7").Value = 705000
$ws.Range("C8").Value = -720000
$ws.Range("C9").Value = 15687500
$ws.Range("C10").Value = 16673619.04761905
$ws.Range("C11").Value = 19320857.14285714
$ws.Range("C12").Value = -2084000
$ws.Range("C13").Value = 68803047.61904763
